$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styles, force Text format so values stay strings (not auto-converted to numbers)
$rng = $ws.Range("D2:G51")
$savedStyle = $rng.Style
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "293.72"
$ws.Range("E2").Value = "0.16%"
$ws.Range("G2").Value = "20"
$ws.Range("D3").Value = "30.99"
$ws.Range("E3").Value = "-0.25%"
$ws.Range("G3").Value = "20"
$ws.Range("D4").Value = "4.907"
$ws.Range("E4").Value = "0.47%"
$ws.Range("G4").Value = "20"
$ws.Range("D5").Value = "0.07334"
$ws.Range("E5").Value = "0.84%"
$ws.Range("G5").Value = "20"
$ws.Range("D6").Value = "2.347"
$ws.Range("E6").Value = "29.75%"
$ws.Range("G6").Value = "20"
$ws.Range("D7").Value = "7.714"
$ws.Range("E7").Value = "0.30%"
$ws.Range("G7").Value = "20"
$ws.Range("D8").Value = "3.738"
$ws.Range("E8").Value = "-0.61%"
$ws.Range("G8").Value = "20"
$ws.Range("D9").Value = "0.9012"
$ws.Range("E9").Value = "-0.34%"
$ws.Range("G9").Value = "20"
$ws.Range("E10").Value = "1.42%"
$ws.Range("G10").Value = "20"
$ws.Range("D11").Value = "0.07892"
$ws.Range("E11").Value = "5.09%"
$ws.Range("G11").Value = "20"
$ws.Range("D12").Value = "0.08147"
$ws.Range("E12").Value = "0.49%"
$ws.Range("G12").Value = "20"
$ws.Range("D13").Value = "0.03099"
$ws.Range("E13").Value = "3.56%"
$ws.Range("G13").Value = "20"
$ws.Range("D14").Value = "0.1007"
$ws.Range("E14").Value = "0.59%"
$ws.Range("G14").Value = "20"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").Value = "0.18%"
$ws.Range("G15").Value = "20"
$ws.Range("D16").Value = "0.005719"
$ws.Range("E16").Value = "-0.75%"
$ws.Range("G16").Value = "20"
$ws.Range("D17").Value = "3.478"
$ws.Range("E17").Value = "0.49%"
$ws.Range("G17").Value = "20"
$ws.Range("E18").Value = "-1.35%"
$ws.Range("G18").Value = "20"
$ws.Range("E19").Value = "0.96%"
$ws.Range("G19").Value = "20"
$ws.Range("D20").Value = "0.1302"
$ws.Range("E20").Value = "-0.35%"
$ws.Range("G20").Value = "20"
$ws.Range("D21").Value = "3.966"
$ws.Range("E21").Value = "-8.80%"
$ws.Range("G21").Value = "20"
$ws.Range("E22").Value = "4.90%"
$ws.Range("G22").Value = "20"
$ws.Range("D23").Value = "0.04537"
$ws.Range("E23").Value = "1.13%"
$ws.Range("G23").Value = "20"
$ws.Range("E24").Value = "-0.96%"
$ws.Range("G24").Value = "20"
$ws.Range("D25").Value = "0.004648"
$ws.Range("E25").Value = "15.00%"
$ws.Range("G25").Value = "20"
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").Value = "3.58%"
$ws.Range("G26").Value = "20"
$ws.Range("D27").Value = "0.0003389"
$ws.Range("G27").Value = "20"
$ws.Range("G28").Value = "20"
$ws.Range("G29").Value = "20"
$ws.Range("G30").Value = "20"
$ws.Range("G31").Value = "20"
$ws.Range("G32").Value = "20"
$ws.Range("G33").Value = "20"
$ws.Range("G34").Value = "20"
$ws.Range("G35").Value = "20"
$ws.Range("G36").Value = "20"
$ws.Range("G37").Value = "20"
$ws.Range("G38").Value = "20"
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "-2.51%"
$ws.Range("G39").Value = "20"
$ws.Range("D40").Value = "0.04447"
$ws.Range("E40").Value = "1.21%"
$ws.Range("G40").Value = "20"
$ws.Range("D41").Value = "0.007322"
$ws.Range("E41").Value = "-1.29%"
$ws.Range("G41").Value = "20"
$ws.Range("D42").Value = "0.1326"
$ws.Range("E42").Value = "0.45%"
$ws.Range("G42").Value = "20"
$ws.Range("D43").Value = "0.008624"
$ws.Range("G43").Value = "20"
$ws.Range("D44").Value = "0.001997"
$ws.Range("E44").Value = "-4.16%"
$ws.Range("G44").Value = "20"
$ws.Range("D45").Value = "0.009432"
$ws.Range("E45").Value = "-7.84%"
$ws.Range("G45").Value = "20"
$ws.Range("D46").Value = "0.00005934"
$ws.Range("E46").Value = "-0.64%"
$ws.Range("G46").Value = "20"
$ws.Range("E47").Value = "-0.36%"
$ws.Range("G47").Value = "20"
$ws.Range("E48").Value = "3.49%"
$ws.Range("G48").Value = "20"
$ws.Range("D49").Value = "0.002894"
$ws.Range("E49").Value = "20.33%"
$ws.Range("G49").Value = "20"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.36%"
$ws.Range("G50").Value = "20"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.36%"
$ws.Range("G51").Value = "20"

# Restore original style so formatting/look is unchanged
$rng.Style = $savedStyle
